$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.275.03"
$ws.Range("E2").Value = "  +5.53%  "
$ws.Range("D3").Value = "2.741.56"
$ws.Range("E3").Value = "  +3.51%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'581.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.65%  "
$ws.Range("D6").Value = "'156.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.67%  "
$ws.Range("D7").Value = "'0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").Value = "'0.610"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.75%  "
$ws.Range("D9").Value = "2.767.78"
$ws.Range("E9").Value = "  +3.85%  "
$ws.Range("D10").Value = "'6.77"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.87%  "
$ws.Range("E11").Value = "  +5.78%  "
$ws.Range("D12").Value = "'0.392"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.77%  "
$ws.Range("E13").Value = "  +2.81%  "
$ws.Range("D14").Value = "3.247.75"
$ws.Range("E14").Value = "  +4.15%  "
$ws.Range("D15").Value = "'27.01"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.67%  "
$ws.Range("D16").Value = "63.865.37"
$ws.Range("E16").Value = "  +4.87%  "
$ws.Range("D17").Value = "'0.0000154"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +7.36%  "
$ws.Range("D18").Value = "2.771.74"
$ws.Range("E18").Value = "  +4.19%  "
$ws.Range("D19").Value = "'12.06"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.62%  "
$ws.Range("D20").Value = "'4.92"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.23%  "
$ws.Range("D21").Value = "'363.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.88%  "
$ws.Range("D22").Value = "'7.02"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.24%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("E24").Value = "  +0.44%  "
$ws.Range("D25").Value = "'66.53"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.81%  "
$ws.Range("E26").Value = "  +5.56%  "
$ws.Range("D27").Value = "'8.59"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.28%  "
$ws.Range("E28").Value = "  +0.40%  "
$ws.Range("D29").Value = "0.0₃0910"
$ws.Range("E29").Value = "  +12.18%  "
$ws.Range("D30").Value = "'2.02"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.83%  "
$ws.Range("D31").Value = "'7.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.01%  "
$ws.Range("D32").Value = "'1.29"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +21.79%  "
$ws.Range("D33").Value = "'173.78"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.12%  "
$ws.Range("D34").Value = "'0.997"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("D35").Value = "'20.61"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.59%  "
$ws.Range("D36").Value = "'4.89"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.85%  "
$ws.Range("D37").Value = "'1.45"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.69%  "
$ws.Range("E38").Value = "  +10.82%  "
$ws.Range("D39").Value = "'1.02"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +13.27%  "
$ws.Range("D40").Value = "'343.34"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.43%  "
$ws.Range("D41").Value = "'4.27"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.68%  "
$ws.Range("D42").Value = "'39.29"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.89%  "
$ws.Range("D43").Value = "'5.86"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +11.39%  "
$ws.Range("D44").Value = "'22.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.58%  "
$ws.Range("D45").Value = "'22.07"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.23%  "
$ws.Range("D46").Value = "'0.0596"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.30%  "
$ws.Range("E47").Value = "  +5.79%  "
$ws.Range("D48").Value = "'138.13"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.07%  "
$ws.Range("D49").Value = "'0.0258"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.90%  "
$ws.Range("D50").Value = "'0.102"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.03%  "
$ws.Range("D51").Value = "'0.998"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.06%  "
